# Generate Report for Handoff
# Appends a new row (for file 8e42f056-ec0a-4c7f-a146-0e3fc7aff847) to each of
# the three worksheets: Overview, zh-cn, de-de. Mirrors the existing row for
# 38f23d60-677f-49ef-a9ec-0d49dee3e928 that's already present in row 2.

$wb = $excel.ActiveWorkbook

$newGuid = "8e42f056-ec0a-4c7f-a146-0e3fc7aff847"
$newHash = "86da594c26610804b305b7406ed3306746770926"

# ---------------------------------------------------------------------
# Sheet "Overview" -> new row 3
# ---------------------------------------------------------------------
$wsOverview = $wb.Worksheets.Item("Overview")

$wsOverview.Hyperlinks.Add(
    $wsOverview.Range("A3"),
    "https://github.com/OpenLocalizationTest/oltest/blob/62a685063ad850ba21bc8b689f56889bb6379ba3/e2e/$newGuid.md",
    "",
    "",
    "$newGuid.md"
) | Out-Null

$wsOverview.Range("B3").Value = "Ready for handoff"
$wsOverview.Range("C3").Value = "Ready for handoff"
$wsOverview.Range("D3").Value = "2016-27-18 10:27:12"

# ---------------------------------------------------------------------
# Sheet "zh-cn" -> new row 3
# ---------------------------------------------------------------------
$wsZhCn = $wb.Worksheets.Item("zh-cn")

$wsZhCn.Hyperlinks.Add(
    $wsZhCn.Range("A3"),
    "https://github.com/OpenLocalizationTest/oltest/blob/62a685063ad850ba21bc8b689f56889bb6379ba3/e2e/$newGuid.md",
    "",
    "",
    "$newGuid.md"
) | Out-Null

$wsZhCn.Hyperlinks.Add(
    $wsZhCn.Range("B3"),
    "https://github.com/OpenLocalizationTest/oltest/blob/62a685063ad850ba21bc8b689f56889bb6379ba3/e2e/$newGuid.md",
    "",
    "",
    ".md"
) | Out-Null

$wsZhCn.Range("C3").Value = "Ready for handoff"

$wsZhCn.Hyperlinks.Add(
    $wsZhCn.Range("D3"),
    "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/6fd77e66b5964f9a7f899c4099aceb11b0e5c6f3/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/ci/ht/$newGuid.$newHash.zh-cn.xlf",
    "",
    "",
    "$newGuid.$newHash.zh-cn.xlf"
) | Out-Null

$wsZhCn.Range("E3").Value = "2016-03-18 10:27:10"
$wsZhCn.Range("E3").NumberFormat = "yyyy-mm-dd HH:mm:ss"

$wsZhCn.Range("H3").Value = "0001-01-01 00:00:00"
$wsZhCn.Range("I3").Value = "Include"

# ---------------------------------------------------------------------
# Sheet "de-de" -> new row 3
# ---------------------------------------------------------------------
$wsDeDe = $wb.Worksheets.Item("de-de")

$wsDeDe.Hyperlinks.Add(
    $wsDeDe.Range("A3"),
    "https://github.com/OpenLocalizationTest/oltest/blob/62a685063ad850ba21bc8b689f56889bb6379ba3/e2e/$newGuid.md",
    "",
    "",
    "$newGuid.md"
) | Out-Null

$wsDeDe.Hyperlinks.Add(
    $wsDeDe.Range("B3"),
    "https://github.com/OpenLocalizationTest/oltest/blob/62a685063ad850ba21bc8b689f56889bb6379ba3/e2e/$newGuid.md",
    "",
    "",
    ".md"
) | Out-Null

$wsDeDe.Range("C3").Value = "Ready for handoff"

$wsDeDe.Hyperlinks.Add(
    $wsDeDe.Range("D3"),
    "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/82d856fdb7423131c55eef42db4e40afd1a845f2/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/ci/ht/$newGuid.$newHash.de-de.xlf",
    "",
    "",
    "$newGuid.$newHash.de-de.xlf"
) | Out-Null

$wsDeDe.Range("E3").Value = "2016-03-18 10:27:12"
$wsDeDe.Range("E3").NumberFormat = "yyyy-mm-dd HH:mm:ss"

$wsDeDe.Range("H3").Value = "0001-01-01 00:00:00"
$wsDeDe.Range("I3").Value = "Include"

Write-Host "Handoff report rows appended."
